# Update gh-pages output data (generated at 456a3b4)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1449
$ws1.Range("F5").Value = 12090
$ws1.Range("F6").Value = 4470
$ws1.Range("F15").Value = 5236
$ws1.Range("F19").Value = 11420
$ws1.Range("F20").Value = 11477

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1449
$ws4.Range("F5").Value = 12090
$ws4.Range("F6").Value = 4470
$ws4.Range("F16").Value = 5236
$ws4.Range("F20").Value = 11420
$ws4.Range("F21").Value = 11477
